$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for Chapter 10 / Chapter 11 rows (rows 11 and 12) ---
$ws.Range("E11").Value = 282
$ws.Range("F11").Value = 335
$ws.Range("H11").Value = 12
$ws.Range("I11").Formula = "=G11/H11"

$ws.Range("E12").Value = 336
$ws.Range("F12").Value = 370
$ws.Range("H12").Value = 6
$ws.Range("I12").Formula = "=G12/H12"

# --- Update summary area (row 1): average now covers G2:G12, drop old L1 estimate ---
$ws.Range("J1").Value = "Ср кол-во стр в день"
$ws.Range("K1").Formula = "=AVERAGE(G2:G12)"
$ws.Range("L1").Clear()

# --- New summary row 2: estimated days remaining ---
$ws.Range("J2").Value = "Оцека кол-ва дней до конца"
$ws.Range("K2").Formula = "=(B15-F12)/K1"
$ws.Range("K2").NumberFormat = "0.0"

# --- New summary row 3: total tomatoes ---
$ws.Range("J3").Value = "Всего помидор"
$ws.Range("K3").Formula = "=SUM(H2:H14)"

# --- Column J needs to fit the new label text ---
$ws.Columns.Item(10).AutoFit() | Out-Null

# --- View adjustments: selection moved from H11 to H13, scrolled back to top ---
$ws.Range("H13").Select() | Out-Null
